$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must remain stored as
# text (matching the workbook's existing inline-string cells), so we
# force Text format before writing, then restore General/Normal so no
# lingering number-format / style is left on the cell.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "244.96"
Set-TextValue "D3" "24.98"
Set-TextValue "D4" "5.123"
Set-TextValue "D5" "0.05658"
Set-TextValue "D6" "6.515"
Set-TextValue "D7" "2.936"
Set-TextValue "D8" "0.8123"
Set-TextValue "D9" "0.8358"
Set-TextValue "D10" "0.1333"
Set-TextValue "D11" "0.06950"
Set-TextValue "D13" "0.09402"
Set-TextValue "D14" "0.001506"
Set-TextValue "D15" "0.0005934"
$ws.Range("E15").Value = "14OneONE"
Set-TextValue "D16" "0.006117"
Set-TextValue "D17" "3.504"
Set-TextValue "D19" "0.3188"
Set-TextValue "D20" "0.03170"
Set-TextValue "D22" "3.746"
Set-TextValue "D25" "0.001235"
Set-TextValue "D26" "0.004261"
Set-TextValue "D27" "0.00009692"
$ws.Range("E27").Value = "26NitroExNTX"
Set-TextValue "D28" "0.0001951"
Set-TextValue "D40" "0.03621"
Set-TextValue "D41" "0.006245"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
Set-TextValue "D42" "0.1050"
Set-TextValue "D43" "0.002720"
Set-TextValue "D44" "0.007387"
Set-TextValue "D45" "0.00005283"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "D47" "0.2198"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue "D48" "0.002285"
Set-TextValue "D49" "0.00002098"
Set-TextValue "D50" "0.0001998"
